$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add lesson #6 ("Methods in Java") as a new row 7 ---

# Start from row 6's formatting (borders/fonts/number formats) so the new
# row matches the rest of the table, then fill in its own values.
$ws.Range("A6:F6").Copy()
$ws.Range("A7:F7").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Enter the YouTube link text first, then the lesson name, so the two new
# shared-string entries are created in the same order as the source file.
$ws.Cells.Item(7, 6).Value = "https://youtu.be/G0B2xyAF3RY"
$ws.Cells.Item(7, 3).Value = "Methods in Java"
$ws.Cells.Item(7, 2).Formula = "=B6+1"
$ws.Cells.Item(7, 4).Value = 2
$ws.Cells.Item(7, 5).Value = 44088

# Row height matches the other single-line rows in the table.
$ws.Rows.Item(7).RowHeight = 15

# Wire up the hyperlink for the new YouTube link, then restore the
# "YouTube link" cell style (Hyperlinks.Add overwrites it with a generic
# hyperlink style) so F7 looks like the other link cells (e.g. F6).
$ws.Hyperlinks.Add($ws.Range("F7"), "https://youtu.be/G0B2xyAF3RY")
$ws.Range("F6").Copy()
$ws.Range("F7").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Column F needs to be a little wider to fit the new link text.
$ws.Columns.Item(6).AutoFit()

# Restore cursor position as left by the author.
$ws.Range("E12").Select()
